$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Products/find" GET method row right before the
#     "Company" section (old row 8 -> new row 9), shifting everything
#     below it down by one row. ---
$ws.Rows("8:8").Insert()

# --- Populate the new row 8 with the Products/find GET method info. ---
$ws.Range("A8").Value = "http://localhost:8080/products/find"
$ws.Range("B8").Value = "GET"
$ws.Range("C8").Value = "desc= ""Descripcion del Producto"""
$ws.Range("D8").Value = "[{`n    ""NAME"": ""BARCEL"",`n    ""DIRECTION"": ""IZCALLI CUAUHTEMOC 3"",`n    ""RFC"": ""JUNA96020""`n}  ]"

# --- Row heights: adjust the new row + all rows that shifted down so the
#     wrapped-text rows keep showing their full content. ---
$ws.Rows(4).RowHeight = 105
$ws.Rows(7).RowHeight = 90
$ws.Rows(8).RowHeight = 75
$ws.Rows(11).RowHeight = 135
$ws.Rows(12).RowHeight = 75
$ws.Rows(13).RowHeight = 150
$ws.Rows(14).RowHeight = 60
$ws.Rows(17).RowHeight = 105
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 105
$ws.Rows(20).RowHeight = 60
$ws.Rows(23).RowHeight = 409.5
$ws.Rows(24).RowHeight = 409.5
$ws.Rows(25).RowHeight = 409.5

# --- Hyperlinks: the row insert above does not shift the hyperlink
#     anchors automatically, so rebuild the whole collection at the new
#     (shifted) cells. Every hyperlinked cell's visible text already IS
#     its target URL, so we can read it back from the cell itself. ---
$ws.Cells.Hyperlinks.Delete()

$linkRows = @(11, 12, 13, 14, 17, 18, 19, 20, 23, 24, 25)
foreach ($r in $linkRows) {
    $cell = $ws.Range("A" + $r)
    $ws.Hyperlinks.Add($cell, $cell.Text)
}
$ws.Hyperlinks.Add($ws.Range("A8"), "http://localhost:8080/products/find")

# Hyperlinks.Add always stamps its own "visited link" style onto the
# cell; re-apply the sheet's normal hyperlink look (copied from A11,
# the Company/new row) on top of every linked A-cell, new one included,
# so they all match the existing method rows.
$ws.Range("A11").Copy()
$styleTargets = @(8, 11, 12, 13, 14, 17, 18, 19, 20, 23, 24, 25)
foreach ($r in $styleTargets) {
    $ws.Range("A" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Move the selection / active cell to D11 (the Company/new response
#     cell) and scroll so row 1 is back in view. ---
$ws.Range("A1").Select()
$ws.Range("D11").Select()
